$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 254.09091
$ws.Range("I33").Value = 194.78947
$ws.Range("K33").Value = 194.78947
$ws.Range("M33").Value = 34.21053000000001
# Row 41
$ws.Range("H41").Value = 71873.71000000001
$ws.Range("I41").Value = 479.23077
$ws.Range("J41").Value = 1000002
$ws.Range("K41").Value = 479.23077
$ws.Range("L41").Value = 1000002
$ws.Range("M41").Value = -39.23077000000001
$ws.Range("N41").Value = -1000882
# Row 138
$ws.Range("H138").Value = 3445.8362
$ws.Range("I138").Value = 2252.0667
$ws.Range("J138").Value = 3835.1086
$ws.Range("K138").Value = 6756.2001
$ws.Range("L138").Value = 11505.3258
$ws.Range("M138").Value = -1616.2001
$ws.Range("N138").Value = -21785.3258

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2454.9333
$ws.Range("I45").Value = 2273.1428
$ws.Range("K45").Value = 2273.1428
$ws.Range("M45").Value = -1896.1428
# Row 74
$ws.Range("H74").Value = 2440.2307
$ws.Range("I74").Value = 1433.75
$ws.Range("K74").Value = 1433.75
$ws.Range("M74").Value = -559.75
# Row 77
$ws.Range("H77").Value = 2440.2307
$ws.Range("I77").Value = 1433.75
$ws.Range("K77").Value = 7168.75
$ws.Range("M77").Value = -2800.75
# Row 102
$ws.Range("H102").Value = 3680.7856
$ws.Range("I102").Value = 3502.3845
$ws.Range("K102").Value = 3502.3845
$ws.Range("M102").Value = -1880.3845
# Row 132
$ws.Range("H132").Value = 1864.2295
$ws.Range("I132").Value = 1602.0182
$ws.Range("K132").Value = 4806.054599999999
$ws.Range("M132").Value = -2276.054599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1597.8334
$ws.Range("I105").Value = 1573.6666
$ws.Range("J105").Value = 1670.3334
$ws.Range("K105").Value = 1573.6666
$ws.Range("L105").Value = 1670.3334
$ws.Range("M105").Value = 173.3334
$ws.Range("N105").Value = -5164.3334
# Row 115
$ws.Range("H115").Value = 144999.5
$ws.Range("J115").Value = 144999.5
$ws.Range("L115").Value = 144999.5
$ws.Range("N115").Value = -148133.5
# Row 134
$ws.Range("H134").Value = 2242.2551
$ws.Range("I134").Value = 1582.1884
$ws.Range("J134").Value = 3812.7585
$ws.Range("K134").Value = 4746.5652
$ws.Range("L134").Value = 11438.2755
$ws.Range("M134").Value = -2211.5652
$ws.Range("N134").Value = -16508.2755

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4018.5667
$ws.Range("I31").Value = 1737.5
$ws.Range("J31").Value = 8580.700000000001
$ws.Range("K31").Value = 1737.5
$ws.Range("L31").Value = 8580.700000000001
$ws.Range("M31").Value = -1442.5
$ws.Range("N31").Value = -9170.700000000001
# Row 34
$ws.Range("H34").Value = 4018.5667
$ws.Range("I34").Value = 1737.5
$ws.Range("J34").Value = 8580.700000000001
$ws.Range("K34").Value = 1737.5
$ws.Range("L34").Value = 8580.700000000001
$ws.Range("M34").Value = -1535.5
$ws.Range("N34").Value = -8984.700000000001
# Row 58
$ws.Range("H58").Value = 1722.9524
$ws.Range("I58").Value = 1088.5483
$ws.Range("K58").Value = 1088.5483
$ws.Range("M58").Value = -885.5482999999999
# Row 107
$ws.Range("H107").Value = 610.5
$ws.Range("I107").Value = 507.375
$ws.Range("K107").Value = 507.375
$ws.Range("M107").Value = 1412.625
# Row 136
$ws.Range("H136").Value = 1722.9524
$ws.Range("I136").Value = 1088.5483
$ws.Range("K136").Value = 3265.6449
$ws.Range("M136").Value = -715.6448999999998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 313.27274
$ws.Range("I10").Value = 334.4
$ws.Range("J10").Value = 102
$ws.Range("K10").Value = 1003.2
$ws.Range("L10").Value = 306
$ws.Range("M10").Value = -864.1999999999999
$ws.Range("N10").Value = -584
# Row 20
$ws.Range("H20").Value = 3831.111
$ws.Range("J20").Value = 4278.75
$ws.Range("L20").Value = 12836.25
$ws.Range("N20").Value = -13290.25
# Row 29
$ws.Range("H29").Value = 1163.7
$ws.Range("J29").Value = 2283.25
$ws.Range("L29").Value = 6849.75
$ws.Range("N29").Value = -7403.75
# Row 122
$ws.Range("H122").Value = 300
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
# Row 133
$ws.Range("H133").Value = 7100.6924

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 999
$ws.Range("J4").Value = 999
$ws.Range("L4").Value = 999
$ws.Range("N4").Value = -1223
# Row 43
$ws.Range("H43").Value = 1749.6666
$ws.Range("I43").Value = 1749.6666
$ws.Range("K43").Value = 1749.6666
$ws.Range("M43").Value = -1598.6666
# Row 57
$ws.Range("H57").Value = 20844.223
$ws.Range("J57").Value = 20844.223
$ws.Range("L57").Value = 20844.223
$ws.Range("N57").Value = -22484.223
# Row 132
$ws.Range("H132").Value = 26325302
$ws.Range("I132").Value = 31256446
$ws.Range("K132").Value = 93769338
$ws.Range("M132").Value = -93766808

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1336.3636
$ws.Range("J46").Value = 1462.625
$ws.Range("L46").Value = 1462.625
$ws.Range("N46").Value = -1838.625
# Row 68
$ws.Range("H68").Value = 3720.3
$ws.Range("I68").Value = 3124.7856
$ws.Range("J68").Value = 5109.8335
$ws.Range("K68").Value = 3124.7856
$ws.Range("L68").Value = 5109.8335
$ws.Range("M68").Value = -2375.7856
$ws.Range("N68").Value = -6607.8335
# Row 71
$ws.Range("H71").Value = 3720.3
$ws.Range("I71").Value = 3124.7856
$ws.Range("J71").Value = 5109.8335
$ws.Range("K71").Value = 15623.928
$ws.Range("L71").Value = 25549.1675
$ws.Range("M71").Value = -11879.928
$ws.Range("N71").Value = -33037.1675
# Row 122
$ws.Range("H122").Value = 5435.885
$ws.Range("I122").Value = 4462.9443
$ws.Range("K122").Value = 13388.8329
$ws.Range("M122").Value = -10938.8329
# Row 132
$ws.Range("H132").Value = 2127.7646
$ws.Range("I132").Value = 1998.1666
$ws.Range("K132").Value = 5994.4998
$ws.Range("M132").Value = -3464.4998
# Row 136
$ws.Range("H136").Value = 2388.6765
$ws.Range("I136").Value = 1391.2
$ws.Range("K136").Value = 4173.6
$ws.Range("M136").Value = -1623.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 125316660
$ws.Range("I4").Value = 1252000
$ws.Range("K4").Value = 1252000
$ws.Range("M4").Value = -1251887
# Row 28
$ws.Range("H28").Value = 5250
$ws.Range("I28").Value = 5000
$ws.Range("J28").Value = 5500
$ws.Range("K28").Value = 5000
$ws.Range("L28").Value = 5500
$ws.Range("M28").Value = -4652
$ws.Range("N28").Value = -6196

# Special case: CUL N122 cell removed entirely (clear contents)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N122").ClearContents()
